$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (codedptmt) holds department codes as text (e.g. "53", "61").
# Force text format first so numeric-looking strings are not auto-converted to numbers.
$ws.Range("C30:C44").NumberFormat = "@"

$ws.Range("A30").Value = "BAILLEUL"
$ws.Range("B30").Value = 43869
$ws.Range("C30").Value = "61"
$ws.Range("D30").Value = 12
$ws.Range("E30").Value = 16.5

$ws.Range("A31").Value = "SAINT PIERRE SUR ERVE"
$ws.Range("B31").Value = 37017
$ws.Range("C31").Value = "53"
$ws.Range("D31").Value = 12
$ws.Range("E31").Value = 12.8

$ws.Range("A32").Value = "RONFEUGERAI"
$ws.Range("B32").Value = 34461
$ws.Range("C32").Value = "61"
$ws.Range("D32").Value = 12
$ws.Range("E32").Value = 10.4

$ws.Range("A33").Value = "ST GERMAIN DU CORBEIS"
$ws.Range("B33").Value = 45704
$ws.Range("C33").Value = "61"
$ws.Range("D33").Value = 11
$ws.Range("E33").Value = 45.5

$ws.Range("A34").Value = "VALFRAMBERT"
$ws.Range("B34").Value = 48801
$ws.Range("C34").Value = "61"
$ws.Range("D34").Value = 11
$ws.Range("E34").Value = 35.5

$ws.Range("A35").Value = "LE BOURGNEUF LA FORET"
$ws.Range("B35").Value = 48399
$ws.Range("C35").Value = "53"
$ws.Range("D35").Value = 11
$ws.Range("E35").Value = 32.5

$ws.Range("A36").Value = "ERNEE"
$ws.Range("B36").Value = 48720
$ws.Range("C36").Value = "53"
$ws.Range("D36").Value = 11
$ws.Range("E36").Value = 21.3

$ws.Range("A37").Value = "COLOMBIERS"
$ws.Range("B37").Value = 47698
$ws.Range("C37").Value = "61"
$ws.Range("D37").Value = 11
$ws.Range("E37").Value = 19.5

$ws.Range("A38").Value = "LONLAY L'ABBAYE"
$ws.Range("B38").Value = 43587
$ws.Range("C38").Value = "61"
$ws.Range("D38").Value = 11
$ws.Range("E38").Value = 19.2

$ws.Range("A39").Value = "SAINT PIERRE DES NIDS"
$ws.Range("B39").Value = 28924
$ws.Range("C39").Value = "53"
$ws.Range("D39").Value = 11
$ws.Range("E39").Value = 19.2

$ws.Range("A40").Value = "AVRILLY"
$ws.Range("B40").Value = 45147
$ws.Range("C40").Value = "61"
$ws.Range("D40").Value = 11
$ws.Range("E40").Value = 17.7

$ws.Range("A41").Value = "BOISSY MAUGIS"
$ws.Range("B41").Value = 43918
$ws.Range("C41").Value = "61"
$ws.Range("D41").Value = 11
$ws.Range("E41").Value = 16.5

$ws.Range("A42").Value = "LA FERRIERE AUX ETANGS"
$ws.Range("B42").Value = 39237
$ws.Range("C42").Value = "61"
$ws.Range("D42").Value = 11
$ws.Range("E42").Value = 15.4

$ws.Range("A43").Value = "ATHIS VAL DE ROUVRE"
$ws.Range("B43").Value = 36300
$ws.Range("C43").Value = "61"
$ws.Range("D43").Value = 11
$ws.Range("E43").Value = 12.8

$ws.Range("A44").Value = "MAYENNE"
$ws.Range("B44").Value = 26478
$ws.Range("C44").Value = "53"
$ws.Range("D44").Value = 12
$ws.Range("E44").Value = 19.2

Write-Host "Applied cyclic shift to rows 30-44"
